$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string (e.g. "215.05") to be stored as
# TEXT, matching the source workbook where every Price/Volume cell is an
# inline/shared string. Plain `.Value = "215.05"` would otherwise be auto-
# detected as the number 215.05 by Excel. Revert the cell style afterwards
# so no visible number-format change is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.895.74"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "1.667.58"
$ws.Range("E3").Value = "  +1.03%  "

Set-TextValue $ws.Range("D5") "215.05"
$ws.Range("E5").Value = "  +0.00%  "

Set-TextValue $ws.Range("D6") "0.521"
$ws.Range("E6").Value = "  +2.37%  "

$ws.Range("E7").Value = "  +0.04%  "

Set-TextValue $ws.Range("D8") "0.0624"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("E9").Value = "  -0.33%  "

Set-TextValue $ws.Range("D10") "20.29"
$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("E11").Value = "  +2.89%  "

$ws.Range("D12").Value = "1.903.00"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").Value = "1.661.34"
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("E14").Value = "  +0.05%  "

Set-TextValue $ws.Range("D15") "0.526"
$ws.Range("E15").Value = "  +1.53%  "

Set-TextValue $ws.Range("D16") "65.67"
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("D17").Value = "26.888.61"
$ws.Range("E17").Value = "  -0.42%  "

Set-TextValue $ws.Range("D18") "234.71"
$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").Value = "0.0₃0732"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("E22").Value = "  -0.44%  "

Set-TextValue $ws.Range("D23") "9.16"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("E24").Value = "  -3.37%  "

Set-TextValue $ws.Range("D25") "146.67"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  -1.12%  "

Set-TextValue $ws.Range("D28") "15.86"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").Value = "1.438.68"
$ws.Range("E33").Value = "  -5.56%  "

$ws.Range("E34").Value = "  +1.61%  "

Set-TextValue $ws.Range("D35") "1.63"
$ws.Range("E35").Value = "  +2.72%  "

$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("E38").Value = "  +2.05%  "

Set-TextValue $ws.Range("D40") "5.73"
$ws.Range("E40").Value = "  -3.67%  "

$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("E42").Value = "  +2.07%  "

Set-TextValue $ws.Range("D43") "0.992"
$ws.Range("E43").Value = "  +8.54%  "

Set-TextValue $ws.Range("D44") "65.93"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").Value = "1.809.98"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("E46").Value = "  +1.38%  "

Set-TextValue $ws.Range("D47") "90.61"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("E49").Value = "  -1.45%  "

Set-TextValue $ws.Range("D50") "0.101"
$ws.Range("E50").Value = "  +3.90%  "

$ws.Range("E51").Value = "  -0.08%  "
